$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 283, shifting the existing rows 283-298 down to 284-299.
$ws.Rows(283).Insert()

# Populate the new row 283 with a fresh weekly price observation
# (same market/product as its neighbours, new date + prices).
$ws.Range("A283").Value = 1
$ws.Range("B283").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C283").Value = "Arica y Parinacota"
$ws.Range("D283").Value = (Get-Date -Year 2022 -Month 9 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E283").Value = 15
$ws.Range("F283").Value = "Fruta"
$ws.Range("G283").Value = 100108
$ws.Range("H283").Value = "Tropicales y subtropicales"
$ws.Range("I283").Value = 100108006
$ws.Range("J283").Value = "Plátano"
$ws.Range("K283").Value = "Sin especificar"
$ws.Range("L283").Value = "Pintón"
$ws.Range("M283").Value = 120
$ws.Range("N283").Value = 25000
$ws.Range("O283").Value = 26000
$ws.Range("P283").Value = 25500
$ws.Range("Q283").Value = "`$/caja 20 kilos"
$ws.Range("R283").Value = "Ecuador"
$ws.Range("S283").Value = 1275
$ws.Range("T283").Value = 20
